$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting the existing rows 23-52 down to 24-53.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new weekly price record.
$ws.Range("A23").Value2 = 9
$ws.Range("B23").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C23").Value2 = "Metropolitana"
$ws.Range("D23").Value2 = 44775
$ws.Range("E23").Value2 = 13
$ws.Range("F23").Value2 = 100112035
$ws.Range("G23").Value2 = "Bruselas (repollito)"
$ws.Range("H23").Value2 = "Sin especificar"
$ws.Range("I23").Value2 = "Primera"
$ws.Range("J23").Value2 = 43
$ws.Range("K23").Value2 = 20000
$ws.Range("L23").Value2 = 20000
$ws.Range("M23").Value2 = 20000
$ws.Range("N23").Value2 = "$/malla 15 kilos"
$ws.Range("O23").Value2 = "Hijuelas"
$ws.Range("P23").Value2 = 1333
$ws.Range("Q23").Value2 = 15
$ws.Range("R23").Value2 = "Hortaliza"
